# feat: add 2022-Q4 data
#
# 1) Insert a new worksheet "2022-Q4" right after "总计" (i.e. as the new
#    second tab), populated with the new quarter's fund-holding table.
# 2) Insert a new row 2 in the "总计" summary sheet for "2022-Q4" and shift
#    the previously existing rows (2022-Q3 .. 2020-Q4) down by one.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: add the new "2022-Q4" sheet, positioned before the current
# second sheet ("2022-Q3"), so the tab order becomes:
#   总计, 2022-Q4, 2022-Q3, 2022-Q2, 2022-Q1, 2021-Q4, 2021-Q3, 2020-Q4
# ---------------------------------------------------------------------
$beforeSheet = $wb.Worksheets.Item(2)
$q4 = $wb.Worksheets.Add($beforeSheet)
$q4.Name = "2022-Q4"
$q4.Outline.SummaryRow = 1
$q4.Outline.SummaryColumn = 1

# Reference sheet to copy cell formatting from (keeps styles identical to
# the rest of the workbook instead of inventing new ones).
$ref = $wb.Worksheets.Item("2022-Q3")

# Header row formatting + text (B1:H1)
$ref.Range("B1:H1").Copy()
$q4.Range("B1:H1").PasteSpecial(-4122)
$q4.Cells.Item(1, 2).Value = "基金代码"
$q4.Cells.Item(1, 3).Value = "基金名称"
$q4.Cells.Item(1, 4).Value = "基金规模"
$q4.Cells.Item(1, 5).Value = "股票总仓位"
$q4.Cells.Item(1, 6).Value = "仓位占比"
$q4.Cells.Item(1, 7).Value = "持有市值(亿元)"
$q4.Cells.Item(1, 8).Value = "仓位排名"

# Row-index column (A2:A3) formatting, matching the other sheets
$ref.Range("A2:A3").Copy()
$q4.Range("A2:A3").PasteSpecial(-4122)
$q4.Cells.Item(2, 1).Value = 0
$q4.Cells.Item(3, 1).Value = 1

# Columns B (fund code) and D:G (numeric-looking figures) are stored as
# text on every other sheet in this workbook, so force text formatting
# before writing - otherwise Excel infers numbers and we lose things like
# leading zeros in fund codes.
$q4.Range("B2:B3").NumberFormat = "@"
$q4.Range("D2:G3").NumberFormat = "@"

$q4.Cells.Item(2, 2).Value = "002446"
$q4.Cells.Item(2, 3).Value = "广发利鑫灵活配置混合A"
$q4.Cells.Item(2, 4).Value = "22.53"
$q4.Cells.Item(2, 5).Value = "73.90"
$q4.Cells.Item(2, 6).Value = "3.01"
$q4.Cells.Item(2, 7).Value = "0.6782"
$q4.Cells.Item(2, 8).Value = 5

$q4.Cells.Item(3, 2).Value = "011172"
$q4.Cells.Item(3, 3).Value = "广发利鑫灵活配置混合C"
$q4.Cells.Item(3, 4).Value = "7.03"
$q4.Cells.Item(3, 5).Value = "73.90"
$q4.Cells.Item(3, 6).Value = "3.01"
$q4.Cells.Item(3, 7).Value = "0.2116"
$q4.Cells.Item(3, 8).Value = 5

# ---------------------------------------------------------------------
# Step 2: insert the 2022-Q4 row into the "总计" (summary) sheet, shifting
# the six existing rows down by one (bottom-up, to not clobber data), then
# write the new row 2.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

for ($r = 7; $r -ge 2; $r--) {
    $dest = $r + 1
    $total.Cells.Item($dest, 2).Value = $total.Cells.Item($r, 2).Value2
    $total.Cells.Item($dest, 3).Value = $total.Cells.Item($r, 3).Value2
    $total.Cells.Item($dest, 4).Value = $total.Cells.Item($r, 4).Value2
    $total.Cells.Item($dest, 1).Value = $r - 1
}

# Row 8 is a brand-new row; copy the index-column (A) style from row 7 so
# it keeps the same bold/centered/bordered look as the rest of column A.
$total.Cells.Item(7, 1).Copy()
$total.Cells.Item(8, 1).PasteSpecial(-4122)

$total.Cells.Item(2, 2).Value = "2022-Q4"
$total.Cells.Item(2, 3).Value = 2
$total.Cells.Item(2, 4).Value = 0.89
$total.Cells.Item(2, 1).Value = 0

# ---------------------------------------------------------------------
# Restore the originally active sheet/tab ("总计") - adding/renaming
# sheets above moves the selection, so put it back where it started.
# ---------------------------------------------------------------------
$total.Activate() | Out-Null
